# Updated cryptos list on Mon Mar 18 07:11:03 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row, and re-ranks three coin pairs whose relative order swapped
# (rows 17/18, 26/27, 42/43), updating Coin name / Link / Price / Volume for
# those rows in full.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (matches the original inlineStr/shared-string cells) and without leaving
# any stray number-format/style behind on the cell.
function Set-TextValue {
    param($Cell, $Value)

    if ($null -eq $Value) { return }

    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Rows whose Coin/Link stay the same - only Price (D) and/or Volume (E)
# text changed.
# ---------------------------------------------------------------------
$rows = @(
    @{ Row = 2;  D = "68.444.97";   E = "  +4.85%  " }
    @{ Row = 3;  D = "3.623.98";    E = "  +4.86%  " }
    @{ Row = 4;  D = $null;         E = "  +0.22%  " }
    @{ Row = 5;  D = "201.37";      E = "  +9.73%  " }
    @{ Row = 6;  D = "584.93";      E = "  +4.08%  " }
    @{ Row = 7;  D = "3.617.27";    E = "  +4.78%  " }
    @{ Row = 8;  D = "0.623";       E = "  +4.21%  " }
    @{ Row = 9;  D = $null;         E = "  -0.18%  " }
    @{ Row = 10; D = "0.686";       E = "  +6.41%  " }
    @{ Row = 11; D = "60.67";       E = "  +18.35%  " }
    @{ Row = 12; D = "0.150";       E = "  +6.07%  " }
    @{ Row = 13; D = "0.0000286";   E = "  +14.01%  " }
    @{ Row = 14; D = "10.15";       E = "  +7.35%  " }
    @{ Row = 15; D = "4.201.71";    E = "  +4.68%  " }
    @{ Row = 16; D = "3.624.29";    E = "  +4.49%  " }
    @{ Row = 19; D = "12.52";       E = "  +6.73%  " }
    @{ Row = 20; D = "68.284.61";   E = "  +5.12%  " }
    @{ Row = 21; D = $null;         E = "  +4.87%  " }
    @{ Row = 22; D = "406.36";      E = "  +6.85%  " }
    @{ Row = 23; D = "13.01";       E = "  +22.14%  " }
    @{ Row = 24; D = "4.30";        E = "  +3.78%  " }
    @{ Row = 25; D = "85.80";       E = "  +3.45%  " }
    @{ Row = 28; D = "12.71";       E = "  +6.34%  " }
    @{ Row = 29; D = "6.14";        E = "  +2.45%  " }
    @{ Row = 30; D = "9.44";        E = "  +10.43%  " }
    @{ Row = 31; D = "7.87";        E = $null }
    @{ Row = 32; D = "31.82";       E = "  +5.47%  " }
    @{ Row = 33; D = "681.51";      E = "  +12.79%  " }
    @{ Row = 34; D = "12.30";       E = "  +4.28%  " }
    @{ Row = 35; D = $null;         E = "  +4.79%  " }
    @{ Row = 36; D = "64.24";       E = "  +2.94%  " }
    @{ Row = 37; D = "42.12";       E = "  +4.38%  " }
    @{ Row = 38; D = "0.418";       E = "  +7.67%  " }
    @{ Row = 39; D = $null;         E = "  -0.15%  " }
    @{ Row = 40; D = "0.0₃0772";    E = "  +6.39%  " }
    @{ Row = 41; D = "3.21";        E = "  +18.36%  " }
    @{ Row = 44; D = "0.998";       E = "  -0.10%  " }
    @{ Row = 45; D = $null;         E = "  +12.57%  " }
    @{ Row = 46; D = "2.90";        E = "  +29.63%  " }
    @{ Row = 47; D = $null;         E = "  +17.02%  " }
    @{ Row = 48; D = "0.0419";      E = "  +6.72%  " }
    @{ Row = 49; D = $null;         E = "  +3.98%  " }
    @{ Row = 50; D = "8.82";        E = "  +8.00%  " }
    @{ Row = 51; D = $null;         E = "  +0.21%  " }
)

foreach ($r in $rows) {
    if ($null -ne $r.D) { Set-TextValue $ws.Cells.Item($r.Row, 4) $r.D }
    if ($null -ne $r.E) { Set-TextValue $ws.Cells.Item($r.Row, 5) $r.E }
}

# ---------------------------------------------------------------------
# Rows that swapped rank order with their neighbour - Coin (B), Link (C),
# Price (D) and Volume (E) all change.
# ---------------------------------------------------------------------
$fullRows = @(
    @{ Row = 17; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "19.32";    E = "  +8.57%  " }
    @{ Row = 18; B = "TRON";      C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx";        D = "0.127";    E = "  +1.37%  " }
    @{ Row = 26; B = "Toncoin";   C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton";          D = "4.01";     E = "  +18.01%  " }
    @{ Row = 27; B = "ImmutableX";C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";       D = "2.94";     E = "  +4.82%  " }
    @{ Row = 42; B = "Maker";     C = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr";        D = "3.198.27"; E = "  +9.16%  " }
    @{ Row = 43; B = "Kaspa";     C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas";             D = "0.135";    E = "  +5.67%  " }
)

foreach ($r in $fullRows) {
    Set-TextValue $ws.Cells.Item($r.Row, 2) $r.B
    Set-TextValue $ws.Cells.Item($r.Row, 3) $r.C
    Set-TextValue $ws.Cells.Item($r.Row, 4) $r.D
    Set-TextValue $ws.Cells.Item($r.Row, 5) $r.E
}
